$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Delete the obsolete sheet
[void]$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()

# Keep the original active/selected sheet
$wb.Worksheets.Item("PAINEIS DARQ").Activate()
